$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2630.2
$ws.Range("I40").Value = 3333.3333
$ws.Range("J40").Value = 2328.8572
$ws.Range("K40").Value = 3333.3333
$ws.Range("L40").Value = 2328.8572
$ws.Range("M40").Value = -3158.3333
$ws.Range("N40").Value = -2678.8572

$ws.Range("H141").Value = 3613.9023
$ws.Range("I141").Value = 2075.8667
$ws.Range("J141").Value = 7808.5454
$ws.Range("K141").Value = 6227.6001
$ws.Range("L141").Value = 23425.6362
$ws.Range("M141").Value = -1047.6001
$ws.Range("N141").Value = -33785.6362

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18699.361
$ws.Range("I32").Value = 4428.6416
$ws.Range("J32").Value = 78458
$ws.Range("K32").Value = 4428.6416
$ws.Range("L32").Value = 78458
$ws.Range("M32").Value = -4141.6416
$ws.Range("N32").Value = -79032

$ws.Range("H45").Value = 724.6667
$ws.Range("I45").Value = 701.875
$ws.Range("J45").Value = 907
$ws.Range("K45").Value = 701.875
$ws.Range("L45").Value = 907
$ws.Range("M45").Value = -324.875
$ws.Range("N45").Value = -1661

$ws.Range("H61").Value = 3958.6191
$ws.Range("I61").Value = 3395.5715
$ws.Range("J61").Value = 6773.857
$ws.Range("K61").Value = 3395.5715
$ws.Range("L61").Value = 6773.857
$ws.Range("M61").Value = -3183.5715
$ws.Range("N61").Value = -7197.857

$ws.Range("H74").Value = 6502.174
$ws.Range("I74").Value = 1008.7143
$ws.Range("J74").Value = 15047.556
$ws.Range("K74").Value = 1008.7143
$ws.Range("L74").Value = 15047.556
$ws.Range("M74").Value = -134.7143
$ws.Range("N74").Value = -16795.556

$ws.Range("H77").Value = 6502.174
$ws.Range("I77").Value = 1008.7143
$ws.Range("J77").Value = 15047.556
$ws.Range("K77").Value = 5043.5715
$ws.Range("L77").Value = 75237.78
$ws.Range("M77").Value = -675.5715
$ws.Range("N77").Value = -83973.78

$ws.Range("H124").Value = 34666.668
$ws.Range("J124").Value = 34666.668
$ws.Range("L124").Value = 34666.668
$ws.Range("N124").Value = -44486.668

$ws.Range("H125").Value = 30216.777
$ws.Range("J125").Value = 30216.777
$ws.Range("L125").Value = 30216.777
$ws.Range("N125").Value = -40056.777

$ws.Range("H136").Value = 3958.6191
$ws.Range("I136").Value = 3395.5715
$ws.Range("J136").Value = 6773.857
$ws.Range("K136").Value = 10186.7145
$ws.Range("L136").Value = 20321.571
$ws.Range("M136").Value = -7636.7145
$ws.Range("N136").Value = -25421.571

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 962.04
$ws.Range("I5").Value = 570.6061
$ws.Range("J5").Value = 1721.8823
$ws.Range("K5").Value = 1711.8183
$ws.Range("L5").Value = 5165.6469
$ws.Range("M5").Value = -1599.8183
$ws.Range("N5").Value = -5389.6469

$ws.Range("H118").Value = 2718.1428
$ws.Range("I118").Value = 1009
$ws.Range("K118").Value = 3027
$ws.Range("M118").Value = -1784

$ws.Range("H122").Value = 1644.7333
$ws.Range("I122").Value = 300.33334
$ws.Range("J122").Value = 2541
$ws.Range("K122").Value = 2703.00006
$ws.Range("L122").Value = 22869
$ws.Range("M122").Value = -253.0000600000003
$ws.Range("N122").Value = -27769

$ws.Range("H135").Value = 962.04
$ws.Range("I135").Value = 570.6061
$ws.Range("J135").Value = 1721.8823
$ws.Range("K135").Value = 5135.4549
$ws.Range("L135").Value = 15496.9407
$ws.Range("M135").Value = -2600.4549
$ws.Range("N135").Value = -20566.9407

$ws.Range("H140").Value = 6608.077
$ws.Range("I140").Value = 8458.27
$ws.Range("J140").Value = 2907.6924
$ws.Range("K140").Value = 25374.81
$ws.Range("L140").Value = 8723.0772
$ws.Range("M140").Value = -20194.81
$ws.Range("N140").Value = -19083.0772

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2855.25
$ws.Range("I80").Value = 2873.3333
$ws.Range("J80").Value = 2801
$ws.Range("K80").Value = 2873.3333
$ws.Range("L80").Value = 2801
$ws.Range("M80").Value = -1875.3333
$ws.Range("N80").Value = -4797

$ws.Range("H83").Value = 2855.25
$ws.Range("I83").Value = 2873.3333
$ws.Range("J83").Value = 2801
$ws.Range("K83").Value = 14366.6665
$ws.Range("L83").Value = 14005
$ws.Range("M83").Value = -9374.666499999999
$ws.Range("N83").Value = -23989

$ws.Range("H107").Value = 1301.4286
$ws.Range("I107").Value = 1301.4286
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1301.4286
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 618.5714
$ws.Range("N107").ClearContents()

$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()

$ws.Range("H122").Value = 2223.077
$ws.Range("I122").Value = 2192
$ws.Range("J122").Value = 2394
$ws.Range("K122").Value = 6576
$ws.Range("L122").Value = 7182
$ws.Range("M122").Value = -4126
$ws.Range("N122").Value = -12082

$ws.Range("H123").Value = 15383.333
$ws.Range("J123").Value = 15383.333
$ws.Range("L123").Value = 15383.333
$ws.Range("N123").Value = -20283.333

$ws.Range("H126").Value = 2820.7144
$ws.Range("J126").Value = 2988.2354
$ws.Range("L126").Value = 8964.706200000001
$ws.Range("N126").Value = -13904.7062

$ws.Range("H132").Value = 2276.923
$ws.Range("I132").Value = 2019.7142
$ws.Range("J132").Value = 3357.2
$ws.Range("K132").Value = 6059.142599999999
$ws.Range("L132").Value = 10071.6
$ws.Range("M132").Value = -3529.142599999999
$ws.Range("N132").Value = -15131.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 632.3929000000001
$ws.Range("I22").Value = 464.72223
$ws.Range("J22").Value = 934.2
$ws.Range("K22").Value = 464.72223
$ws.Range("L22").Value = 934.2
$ws.Range("M22").Value = -169.72223
$ws.Range("N22").Value = -1524.2

$ws.Range("H27").Value = 632.3929000000001
$ws.Range("I27").Value = 464.72223
$ws.Range("J27").Value = 934.2
$ws.Range("K27").Value = 464.72223
$ws.Range("L27").Value = 934.2
$ws.Range("M27").Value = -357.72223
$ws.Range("N27").Value = -1148.2

$ws.Range("H68").Value = 2004
$ws.Range("I68").Value = 1891.4286
$ws.Range("J68").Value = 2266.6667
$ws.Range("K68").Value = 1891.4286
$ws.Range("L68").Value = 2266.6667
$ws.Range("M68").Value = -1142.4286
$ws.Range("N68").Value = -3764.6667

$ws.Range("H71").Value = 2004
$ws.Range("I71").Value = 1891.4286
$ws.Range("J71").Value = 2266.6667
$ws.Range("K71").Value = 9457.143
$ws.Range("L71").Value = 11333.3335
$ws.Range("M71").Value = -5713.143
$ws.Range("N71").Value = -18821.3335

$ws.Range("H82").Value = 2213.9285
$ws.Range("I82").Value = 2441.7144
$ws.Range("J82").Value = 1986.1428
$ws.Range("K82").Value = 2441.7144
$ws.Range("L82").Value = 1986.1428
$ws.Range("M82").Value = -2080.7144
$ws.Range("N82").Value = -2708.1428

$ws.Range("H85").Value = 2213.9285
$ws.Range("I85").Value = 2441.7144
$ws.Range("J85").Value = 1986.1428
$ws.Range("K85").Value = 2441.7144
$ws.Range("L85").Value = 1986.1428
$ws.Range("M85").Value = -1193.7144
$ws.Range("N85").Value = -4482.1428

$ws.Range("H93").Value = 1888.4783
$ws.Range("I93").Value = 1792.7142
$ws.Range("J93").Value = 2037.4445
$ws.Range("K93").Value = 1792.7142
$ws.Range("L93").Value = 2037.4445
$ws.Range("M93").Value = -544.7141999999999
$ws.Range("N93").Value = -4533.4445

$ws.Range("H136").Value = 3014.861
$ws.Range("I136").Value = 1701
$ws.Range("J136").Value = 9584.166999999999
$ws.Range("K136").Value = 5103
$ws.Range("L136").Value = 28752.501
$ws.Range("M136").Value = -2553
$ws.Range("N136").Value = -33852.501

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1743.9072
$ws.Range("I132").Value = 1812.5325
$ws.Range("J132").Value = 1479.7
$ws.Range("K132").Value = 5437.5975
$ws.Range("L132").Value = 4439.1
$ws.Range("M132").Value = -2907.5975
$ws.Range("N132").Value = -9499.1
